# Detect inactive users in shared game play.
# Adds three new template rows (075, 076, 077) to the "Templates" sheet,
# describing the new "shared play - preparing" state and the
# "confirm cancel dialog" that lets a user cancel out of a shared-play
# session while waiting on other (possibly inactive) players.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")

# --- New row 72: template 075 - shared play - preparing -------------------
$ws.Range("A72").Value = "075"
$ws.Range("B72").Value = 1273
$ws.Range("C72").Value = 838
$ws.Range("D72").Value = 1518
$ws.Range("E72").Value = 903
$ws.Range("F72").Value = "075"
$ws.Range("G72").Value = "shared play - preparing "

# --- New row 73: template 076 - confirm cancel dialog - title -------------
$ws.Range("A73").Value = "076"
$ws.Range("B73").Value = 743
$ws.Range("C73").Value = 321
$ws.Range("D73").Value = 1103
$ws.Range("E73").Value = 373
$ws.Range("F73").Value = "076"
$ws.Range("G73").Value = "confirm cancel dialog - title"

# --- New row 74: template 077 - confirm cancel dialog - ok button ---------
$ws.Range("A74").Value = "077"
$ws.Range("B74").Value = 1080
$ws.Range("C74").Value = 678
$ws.Range("D74").Value = 1178
$ws.Range("E74").Value = 725
# NOTE: matches the existing authoring quirk seen elsewhere in this column
# (e.g. row 69), where File Name (F) reuses the "037" label instead of the
# Template Number of its own row.
$ws.Range("F74").Value = "037"
$ws.Range("G74").Value = "confirm cancel dialog - ok button"

# --- Refresh the view: scroll down to the newly added rows and select them
$ws.Activate()
$ws.Range("C66").Select()
$ws.Range("A1:G74").Select()
